$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 1 ---
$ws.Range("A1").Value = 158064
$ws.Range("B1").Value = "Exercícios do livro de matemática"
$ws.Range("C1").Value = 78
$ws.Range("D1").Value = 51

# Preserve / reapply the distinctive formatting of C1 (Arial 10, wrap text, General format)
$ws.Range("C1").WrapText = $true
$ws.Range("C1").Font.Name = "Arial"
$ws.Range("C1").Font.Size = 10

# --- Update row 2 ---
$ws.Range("A2").Value = 163070
$ws.Range("B2").Value = "Exercício01_atualizado"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 5

# --- Add row 3 (new data row) ---
$ws.Range("A3").Value = 163073
$ws.Range("B3").Value = "Exercício04_atualizado"
$ws.Range("C3").Value = 40
$ws.Range("D3").Value = 20

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 41.27
$ws.StandardWidth = 8.7578125

# --- Selection ---
$ws.Range("E5").Select() | Out-Null
